# "work on localconf passwords"
#
# The two ALTER-TABLE note cells (row 6 / row 17) move from column A to
# column B, lose their highlighted "note" style, and have their text
# trimmed down to just the trailing UNIQUE (...) clause. Column B also
# gets wider to accommodate the new values column header, and the
# worksheet's selection moves from A18 to B20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: storagepass constraint note -> shortened, moved A6 -> B6 ---
$ws.Rows(6).ClearFormats()
$ws.Range("A6").ClearContents()
$ws.Range("B6").Value = "UNIQUE (vfs , pool)"

# --- Row 17: localconf constraint note -> shortened, moved A17 -> B17 ---
$ws.Rows(17).ClearFormats()
$ws.Range("A17").ClearContents()
$ws.Range("B17").Value = "UNIQUE (service , section , key)"

# --- Column B widened to fit the longer "values" entries ---
$ws.Columns(2).ColumnWidth = 29.666666666666668

# --- Selection moves to B20 ---
$ws.Range("B20").Select()
